$wb = $excel.ActiveWorkbook

# 1) Status text "Ready for handoff" -> "In Translation" everywhere it appears
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# 2) Narrow the "Status" columns (Overview columns E & F, and column C on the
#    language sheets) from their previous wide setting (~17.22 chars) to the
#    new, narrower one (~13.41 chars) to match the shorter status text.
#    ColumnWidth assignments are snapped to the host's internal pixel grid,
#    so 12.5 is the input that lands closest to the target width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
